# Apply cryptos list update (price/volume refresh) as per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.605.21'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '2.388.46'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.29'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.70'
$ws.Range('E6').Value = '  +4.09%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '2.387.51'
$ws.Range('E9').Value = '  +2.25%  '
$ws.Range('E10').Value = '  +3.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.157'
$ws.Range('E11').Value = '  +1.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.35'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  +3.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.86'
$ws.Range('E14').Value = '  +5.89%  '
$ws.Range('E15').Value = '  +8.92%  '
$ws.Range('D16').Value = '2.821.26'
$ws.Range('E16').Value = '  +2.37%  '
$ws.Range('D17').Value = '61.574.34'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').Value = '2.389.28'
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.10'
$ws.Range('E19').Value = '  +5.80%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.53'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.16'
$ws.Range('E21').Value = '  +2.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.66'
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.48'
$ws.Range('E24').Value = '  +2.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.74'
$ws.Range('E25').Value = '  -5.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.19'
$ws.Range('E26').Value = '  +7.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '553.36'
$ws.Range('E27').Value = '  +11.03%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.28'
$ws.Range('E29').Value = '  +4.55%  '
$ws.Range('D30').Value = '2.473.75'
$ws.Range('D31').Value = '0.0₃0917'
$ws.Range('E31').Value = '  +3.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  +1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.148'
$ws.Range('E33').Value = '  +2.73%  '
$ws.Range('E34').Value = '  +3.59%  '
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.72'
$ws.Range('E36').Value = '  +9.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.93'
$ws.Range('E38').Value = '  +7.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.73'
$ws.Range('E39').Value = '  +3.18%  '
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.54'
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '145.47'
$ws.Range('E42').Value = '  +5.83%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  +7.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '148.28'
$ws.Range('E45').Value = '  +5.26%  '
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0528'
$ws.Range('E47').Value = '  +3.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.18'
$ws.Range('E48').Value = '  +4.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.584'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0906'
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('E51').Value = '  +1.56%  '
